$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: 展览 (Exhibitions) -- bump "want-to-go" counts (column F)
# ---------------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 497
$wsExhibit.Range("F3").Value  = 1518
$wsExhibit.Range("F4").Value  = 792
$wsExhibit.Range("F7").Value  = 1102
$wsExhibit.Range("F8").Value  = 689
$wsExhibit.Range("F10").Value = 1347
$wsExhibit.Range("F12").Value = 1009
$wsExhibit.Range("F16").Value = 41
$wsExhibit.Range("F20").Value = 529

# ---------------------------------------------------------------------------
# Sheet: 演出 (Performances) -- bump "want-to-go" counts (column F)
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value  = 131
$wsShow.Range("F9").Value  = 581
$wsShow.Range("F10").Value = 52

# ---------------------------------------------------------------------------
# Sheet: 全部类型 (All types) -- bump "want-to-go" counts (column F)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 497
$wsAll.Range("F4").Value  = 1518
$wsAll.Range("F10").Value = 1102
$wsAll.Range("F11").Value = 689
$wsAll.Range("F13").Value = 1347
$wsAll.Range("F15").Value = 1009
$wsAll.Range("F19").Value = 41
$wsAll.Range("F26").Value = 131
$wsAll.Range("F27").Value = 529

# Row 30 (index 29) becomes what used to be row 31's event
# (广州·代号鸢only2.0 duplicate in row 29 is left untouched)
$wsAll.Range("C30").Value = "广州·原神X星穹铁道X绝区零ONLY"
$wsAll.Range("D30").Value = "洛浦街夏滘西环路1号(厦滘地铁站A口步行290米) 厦喾岭南电商园会展中心"
$wsAll.Range("E30").Value = "2024.03.16 10:00-03.16 17:00"
$wsAll.Range("F30").Value = 222
$wsAll.Range("G30").Value = 60
$wsAll.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=80715"
$wsAll.Range("I30").Value = "//i0.hdslb.com/bfs/openplatform/202401/Lt6ZYvA41704878219924.jpeg"

# Row 31 becomes what used to be row 32's event
# (force literal text so Excel does not reinterpret the date-like string
# as a real date value/format, matching the original inlineStr storage)
$wsAll.Range("B31").NumberFormat = "@"
$wsAll.Range("B31").Value = "2024-03-17"
$wsAll.Range("B31").Style = "Normal"
$wsAll.Range("C31").Value = "广州·三月的幻想演唱会2024「飞越蓝色时刻」"
$wsAll.Range("D31").Value = "恩宁路265号三层、四层自编01 MAO Livehouse广州(永庆坊店)"
$wsAll.Range("E31").Value = "2024.03.17 19:00-03.17 20:30"
$wsAll.Range("F31").Value = 59
$wsAll.Range("G31").Value = 380
$wsAll.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=80870"
$wsAll.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202401/8WBT7H6W1705376580145.png"

# Row 32 becomes what used to be row 33's event
$wsAll.Range("B32").NumberFormat = "@"
$wsAll.Range("B32").Value = "2024-03-23"
$wsAll.Range("B32").Style = "Normal"
$wsAll.Range("C32").Value = "广州·排球少年ONLY"
$wsAll.Range("D32").Value = "机场路1399号广州百信广场二期 李宁运动中心"
$wsAll.Range("E32").Value = "2024.03.23 10:00-03.23 17:00"
$wsAll.Range("F32").Value = 162
$wsAll.Range("G32").Value = 60
$wsAll.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=80716"
$wsAll.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202401/IFLvYmxx1704879325152.jpeg"

# Row 33 becomes a brand-new event (春卷饭十周年专场演出)
$wsAll.Range("C33").Value = "广州·春卷饭 十周年  2024  专场演出"
$wsAll.Range("D33").Value = "革新路124号太古仓码头54汇5号仓 太空间Livehouse"
$wsAll.Range("E33").Value = "2024.03.23 20:00-03.23 22:00"
$wsAll.Range("F33").Value = 581
$wsAll.Range("G33").Value = "已售罄"
$wsAll.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=81186"
$wsAll.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202401/ho9rIMg21705894649801.jpeg"

# Rows 34/35 (KANAKO ITO&AYANE duplicate entries) -- bump want-to-go count
$wsAll.Range("F34").Value = 52
$wsAll.Range("F35").Value = 52
